# Union base de datos DUQUE al dataframe y graficos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# Copy the formatting (date format, wrap text, numeric style) from the last
# existing row (96) down onto the three new rows so the new cells pick up
# the same styles already used by the table (s="18"/"2"/"9").
$ws.Range("A96:C96").Copy()
$ws.Range("A97:C99").PasteSpecial(-4122)

# Row 97: date 44786 -> 2022-08-13
$ws.Cells.Item(97, 1).Value = 44786
$ws.Cells.Item(97, 2).Value = "Nuevos graficos articulo ICC, lectura de articulo profe mando"
$ws.Cells.Item(97, 3).Value = 3

# Row 98: date 44788 -> 2022-08-15
$ws.Cells.Item(98, 1).Value = 44788
$ws.Cells.Item(98, 2).Value = "Cedulas de la base de datos de duque, asesoria de ASC con los estudiantes."
$ws.Cells.Item(98, 3).Value = 4

# Row 99: date 44793 -> 2022-08-20
$ws.Cells.Item(99, 1).Value = 44793
$ws.Cells.Item(99, 2).Value = "Reunion con vero, lectura de articulos"
$ws.Cells.Item(99, 3).Value = 6

# Rows 97 and 98 wrap onto two lines like row 96, so they use the taller
# row height; row 99's text fits on one line and keeps the default height.
$ws.Rows.Item(97).RowHeight = 27.6
$ws.Rows.Item(98).RowHeight = 27.6

# Scroll the view down and select the new last cell, matching how the
# sheet was left after the new entries were appended.
$ws.Application.ActiveWindow.ScrollRow = 88
$ws.Range("C99").Select()
